# Trade #67 closed at 2026-02-18 00:24:14 - unknown UNKNOWN +0.000%
#
# This script mirrors the commit's effect on the workbook:
#   - Trade #95 (a MarketMaking trade) moves from OPEN -> CLOSED (early_exit)
#     on both the "All Trades" sheet and the strategy-specific "MarketMaking"
#     sheet.
#   - Summary / Strategy Status roll-up counters are refreshed accordingly.
#   - A brand-new OPEN trade (#124) is appended to both the "All Trades" and
#     "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 95        # Total Trades
$summary.Range("B9").Value = 47.37     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 6 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 31         # Trades
$status.Range("G6").Value = 48.39      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - close trade #95 (row 96)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G96").Value = 0.009423000000000001   # Exit Price
$allTrades.Range("H96").Value = "CLOSED"                # Status
$allTrades.Range("I96").Value = -5.7671                 # P&L %
$allTrades.Range("J96").Value = -0                      # P&L $
$allTrades.Range("K96").Value = 99.41                    # Capital After
$allTrades.Range("L96").Value = "early_exit"             # Exit Reason
$allTrades.Range("M96").Value = 0.12                      # Duration (min)

# Append new trade #124 (row 125) to All Trades
$allTrades.Range("A125").Value = 124
$allTrades.Range("B125").NumberFormat = "@"
$allTrades.Range("B125").Value = "2026-02-18"
$allTrades.Range("C125").NumberFormat = "@"
$allTrades.Range("C125").Value = "00:24:08"
$allTrades.Range("D125").Value = "MarketMaking"
$allTrades.Range("E125").Value = "UP"
$allTrades.Range("F125").Value = 0.01
$allTrades.Range("H125").Value = "OPEN"
$allTrades.Range("I125").Value = 0
$allTrades.Range("J125").Value = 0
$allTrades.Range("K125").Value = 99.410254715139
$allTrades.Range("M125").Value = 0
$allTrades.Range("N125").Value = 0
$allTrades.Range("O125").Value = 0
$allTrades.Range("P125").Value = 0.6
$allTrades.Range("Q125").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet - close trade #95 (row 32)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G32").Value = 0.009423000000000001   # Exit Price
$mm.Range("H32").Value = "CLOSED"                # Status
$mm.Range("I32").Value = -5.7671                 # P&L %
$mm.Range("J32").Value = -0                      # P&L $
$mm.Range("K32").Value = 99.41                    # Capital After
$mm.Range("P32").Value = "early_exit"             # Exit Reason
$mm.Range("Q32").Value = 0.12                      # Duration (min)

# Append new trade #124 (row 45) to MarketMaking
$mm.Range("A45").Value = 124
$mm.Range("B45").NumberFormat = "@"
$mm.Range("B45").Value = "2026-02-18"
$mm.Range("C45").NumberFormat = "@"
$mm.Range("C45").Value = "00:24:08"
$mm.Range("D45").Value = "MarketMaking"
$mm.Range("E45").Value = "UP"
$mm.Range("F45").Value = 0.01
$mm.Range("H45").Value = "OPEN"
$mm.Range("I45").Value = 0
$mm.Range("J45").Value = 0
$mm.Range("K45").Value = 99.410254715139
$mm.Range("L45").Value = 0
$mm.Range("M45").Value = 0
$mm.Range("N45").Value = 0.6
$mm.Range("O45").Value = "Normal spread capture: 198 bps"
$mm.Range("Q45").Value = 0

Write-Host "Applied trade #95 close + trade #124 open edits"
